$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the previously-empty runtime measurements (columns B:D) that were
#    causing the #DIV/0! errors in column G.
# ---------------------------------------------------------------------------
$ws.Range("B4").Value  = 0.221802933
$ws.Range("B5").Value  = 0.204817407
$ws.Range("B6").Value  = 0.199692448
$ws.Range("B7").Value  = 0.167936961
$ws.Range("B8").Value  = 0.202551462
$ws.Range("B9").Value  = 0.211043682
$ws.Range("B11").Value = 0.625955684
$ws.Range("C11").Value = 0.592214591
$ws.Range("D11").Value = 0.588886076
$ws.Range("B13").Value = 0.565636734
$ws.Range("C13").Value = 0.603483458

# ---------------------------------------------------------------------------
# 2. Add the "Speedup" column (H) next to the existing "Durchschnitt" (G)
#    column: H1 header + H4 anchor formula + shared H5:H13 formula (skipping
#    the blank rows 10/12 that only hold the section labels).
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "Speedup"
$ws.Range("H4").Formula = "=G$4/G4"
$ws.Range("H5:H13").Formula = "=G`$4/G5"
$ws.Range("H10").ClearContents()
$ws.Range("H12").ClearContents()

# ---------------------------------------------------------------------------
# 3. Build the second ("bad") speedup table in columns J:Q, mirroring the
#    A:H layout exactly (J<->A, K<->B, P<->G, Q<->H).
# ---------------------------------------------------------------------------
# Row 1 headers (copy format from B1:F1 / G1 so the merged style is reused,
# then restore the "Durchlauf" text that PasteSpecial(formats) strips out --
# matching the source B1 which also carries the text under the merge).
$ws.Range("B1:F1").Copy() | Out-Null
$ws.Range("K1:O1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("K1:O1").Merge()
$ws.Range("K1").Value = "Durchlauf"
$ws.Range("P1").Value = "Durchschnitt"
$ws.Range("Q1").Value = "Speedup"

# Row 2 (run index headers 1..5)
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 5

# Section labels mirrored into column J
$ws.Range("J3").Value  = "1 Place"
$ws.Range("J10").Value = "2 Places"
$ws.Range("J12").Value = "4 Places"

# Place-count column (J) + raw measurements (K)
$ws.Range("J4").Value  = 1
$ws.Range("K4").Value  = 23.713405121
$ws.Range("J5").Value  = 2
$ws.Range("J6").Value  = 4
$ws.Range("K6").Value  = 223.026686391
$ws.Range("J7").Value  = 8
$ws.Range("K7").Value  = 379.691890248
$ws.Range("J8").Value  = 16
$ws.Range("K8").Value  = 389.670922313
$ws.Range("J9").Value  = 32
$ws.Range("K9").Value  = 397.110874808
$ws.Range("J11").Value = 32
$ws.Range("K11").Value = 199.253641429
$ws.Range("J13").Value = 32
$ws.Range("K13").Value = 103.400613564

# Durchschnitt (P) + Speedup (Q) formulas, mirroring G/H.
# NOTE: written *before* the K-column NumberFormat below, otherwise the new
# formula cells would inherit the "#,##0" number format from their K
# precedent cell.
$ws.Range("P4").Formula = "=SUM(K4:O4)/COUNTA(K4:O4)"
$ws.Range("P5:P13").Formula = "=SUM(K5:O5)/COUNTA(K5:O5)"
$ws.Range("P10").ClearContents()
$ws.Range("P12").ClearContents()

$ws.Range("Q4").Formula = "=P$4/P4"
$ws.Range("Q5:Q13").Formula = "=P`$4/P5"
$ws.Range("Q10").ClearContents()
$ws.Range("Q12").ClearContents()

# Parameter legend row (15) mirrored into J:P
$ws.Range("J15").Value = "Parameter:"
$ws.Range("K15").Value = "n"
$ws.Range("L15").Value = "m"
$ws.Range("M15").Value = "seedA"
$ws.Range("N15").Value = "d"
$ws.Range("O15").Value = "i"
$ws.Range("P15").Value = "verbose"

# ---------------------------------------------------------------------------
# 4. New example-parameters row (16) for both tables.
# ---------------------------------------------------------------------------
$ws.Range("B16").Value = 4
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 1

$ws.Range("K16").Value = 100
$ws.Range("L16").Value = 10
$ws.Range("M16").Value = 42
$ws.Range("N16").Value = 3
$ws.Range("O16").Value = 10
$ws.Range("P16").Value = 0

# ---------------------------------------------------------------------------
# 5. Apply the thousands-separator number format to the raw K measurements
#    only -- done last so it cannot leak onto the P/Q formula cells that
#    reference column K.
# ---------------------------------------------------------------------------
foreach ($r in @(4,6,7,8,9,11,13)) {
    $ws.Range("K$r").NumberFormat = "#,##0"
}

# ---------------------------------------------------------------------------
# 6. Selection mirrors the author's last selection when they saved.
# ---------------------------------------------------------------------------
$ws.Range("K5").Select() | Out-Null

Write-Output "edit applied"
